$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 500
$ws.Range("J97").Value = 500
$ws.Range("L97").Value = 1500
$ws.Range("N97").Value = -2492
$ws.Range("H103").Value = 1127.6
$ws.Range("I103").Value = 669.3333
$ws.Range("J103").Value = 1815
$ws.Range("K103").Value = 2007.9999
$ws.Range("L103").Value = 5445
$ws.Range("M103").Value = -1421.9999
$ws.Range("N103").Value = -6617
$ws.Range("H106").Value = 4610.636
$ws.Range("I106").Value = 2838.125
$ws.Range("J106").Value = 9337.333000000001
$ws.Range("K106").Value = 2838.125
$ws.Range("L106").Value = 9337.333000000001
$ws.Range("M106").Value = -2207.125
$ws.Range("N106").Value = -10599.333
$ws.Range("H112").Value = 1335.8857
$ws.Range("J112").Value = 1392.4062
$ws.Range("L112").Value = 4177.2186
$ws.Range("N112").Value = -6393.2186
$ws.Range("H123").Value = 44657.5
$ws.Range("J123").Value = 44657.5
$ws.Range("L123").Value = 44657.5
$ws.Range("N123").Value = -54457.5
$ws.Range("H128").Value = 25660
$ws.Range("J128").Value = 25660
$ws.Range("L128").Value = 25660
$ws.Range("N128").Value = -35620
$ws.Range("H129").Value = 1100.293
$ws.Range("I129").Value = 565.6667
$ws.Range("J129").Value = 1117
$ws.Range("K129").Value = 1697.0001
$ws.Range("L129").Value = 3351
$ws.Range("M129").Value = 3302.9999
$ws.Range("N129").Value = -13351
$ws.Range("H137").Value = 2750.0833
$ws.Range("I137").Value = 1242.1666
$ws.Range("J137").Value = 4258
$ws.Range("K137").Value = 3726.4998
$ws.Range("L137").Value = 12774
$ws.Range("M137").Value = -1176.4998
$ws.Range("N137").Value = -17874

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 4800
$ws.Range("I31").Value = 4800
$ws.Range("K31").Value = 4800
$ws.Range("M31").Value = -4506
$ws.Range("H74").Value = 1306.5853
$ws.Range("I74").Value = 1246.6471
$ws.Range("K74").Value = 1246.6471
$ws.Range("M74").Value = -372.6470999999999
$ws.Range("H77").Value = 1306.5853
$ws.Range("I77").Value = 1246.6471
$ws.Range("K77").Value = 6233.2355
$ws.Range("M77").Value = -1865.2355
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0
$ws.Range("H141").Value = 72809.664
$ws.Range("J141").Value = 72809.664
$ws.Range("L141").Value = 72809.664
$ws.Range("N141").Value = -83169.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2466.66
$ws.Range("I105").Value = 1280.591
$ws.Range("J105").Value = 2801.1924
$ws.Range("K105").Value = 1280.591
$ws.Range("L105").Value = 2801.1924
$ws.Range("M105").Value = 466.4090000000001
$ws.Range("N105").Value = -6295.1924
$ws.Range("H134").Value = 1430.7556
$ws.Range("I134").Value = 1341.9722
$ws.Range("K134").Value = 4025.9166
$ws.Range("M134").Value = -1490.9166
$ws.Range("H137").Value = 8000
$ws.Range("I137").Value = 8000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 8000
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0
$ws.Range("M137").Value = -2900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4802.5
$ws.Range("I16").Value = 1605.5
$ws.Range("J16").Value = 7999.5
$ws.Range("K16").Value = 1605.5
$ws.Range("L16").Value = 7999.5
$ws.Range("M16").Value = -1318.5
$ws.Range("N16").Value = -8573.5
$ws.Range("H19").Value = 2381.4285
$ws.Range("I19").Value = 145
$ws.Range("J19").Value = 5363.3335
$ws.Range("K19").Value = 145
$ws.Range("L19").Value = 5363.3335
$ws.Range("M19").Value = 25
$ws.Range("N19").Value = -5703.3335
$ws.Range("H24").Value = 2381.4285
$ws.Range("I24").Value = 145
$ws.Range("J24").Value = 5363.3335
$ws.Range("K24").Value = 145
$ws.Range("L24").Value = 5363.3335
$ws.Range("M24").Value = 25
$ws.Range("N24").Value = -5703.3335
$ws.Range("H113").Value = 4802.5
$ws.Range("I113").Value = 1605.5
$ws.Range("J113").Value = 7999.5
$ws.Range("K113").Value = 1605.5
$ws.Range("L113").Value = 7999.5
$ws.Range("M113").Value = 564.5
$ws.Range("N113").Value = -12339.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 4619.25
$ws.Range("I110").Value = 3370.8
$ws.Range("J110").Value = 6700
$ws.Range("K110").Value = 10112.4
$ws.Range("L110").Value = 20100
$ws.Range("M110").Value = -6022.400000000001
$ws.Range("N110").Value = -28280
$ws.Range("H121").Value = 611.1818
$ws.Range("I121").Value = 576.6667
$ws.Range("J121").Value = 766.5
$ws.Range("K121").Value = 1730.0001
$ws.Range("L121").Value = 2299.5
$ws.Range("M121").Value = -420.0001
$ws.Range("N121").Value = -4919.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3350334.8
$ws.Range("I14").Value = 5000502
$ws.Range("K14").Value = 5000502
$ws.Range("M14").Value = -5000334
$ws.Range("H70").Value = 5690.2285
$ws.Range("I70").Value = 5919.2856
$ws.Range("J70").Value = 4774
$ws.Range("K70").Value = 5919.2856
$ws.Range("L70").Value = 4774
$ws.Range("M70").Value = -5649.2856
$ws.Range("N70").Value = -5314
$ws.Range("H73").Value = 5690.2285
$ws.Range("I73").Value = 5919.2856
$ws.Range("J73").Value = 4774
$ws.Range("K73").Value = 5919.2856
$ws.Range("L73").Value = 4774
$ws.Range("M73").Value = -4983.2856
$ws.Range("N73").Value = -6646
$ws.Range("H97").Value = 3124.111
$ws.Range("I97").Value = 3765.6667
$ws.Range("J97").Value = 2803.3333
$ws.Range("K97").Value = 3765.6667
$ws.Range("L97").Value = 2803.3333
$ws.Range("M97").Value = -3269.6667
$ws.Range("N97").Value = -3795.3333
$ws.Range("H132").Value = 2398
$ws.Range("I132").Value = 2105.2122
$ws.Range("J132").Value = 3088.1428
$ws.Range("K132").Value = 6315.6366
$ws.Range("L132").Value = 9264.428400000001
$ws.Range("M132").Value = -3785.6366
$ws.Range("N132").Value = -14324.4284
$ws.Range("H138").Value = 30429
$ws.Range("J138").Value = 30429
$ws.Range("L138").Value = 30429
$ws.Range("N138").Value = -40709
$ws.Range("H140").Value = 46500
$ws.Range("I140").Value = 40000
$ws.Range("J140").Value = 53000
$ws.Range("K140").Value = 40000
$ws.Range("L140").Value = 53000
$ws.Range("M140").Value = -34820
$ws.Range("N140").Value = -63360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3500
$ws.Range("J20").Value = 3500
$ws.Range("L20").Value = 3500
$ws.Range("N20").Value = -3952
$ws.Range("H36").Value = 74140
$ws.Range("J36").Value = 74140
$ws.Range("L36").Value = 74140
$ws.Range("N36").Value = -75264
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0
$ws.Range("H138").Value = 47294
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 47294
$ws.Range("K138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("M138").Value = 47294
$ws.Range("N138").Value = -57574

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 420.125
$ws.Range("I113").Value = 420.125
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1260.375
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 909.625
$ws.Range("H122").Value = 1991.36
$ws.Range("I122").Value = 1332.6364
$ws.Range("J122").Value = 2508.9285
$ws.Range("K122").Value = 3997.9092
$ws.Range("L122").Value = 7526.7855
$ws.Range("M122").Value = -1547.9092
$ws.Range("N122").Value = -12426.7855
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0
$ws.Range("H132").Value = 1084.2963
$ws.Range("I132").Value = 686.561
$ws.Range("K132").Value = 2059.683
$ws.Range("M132").Value = 470.317
$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -39920
